# The commit swaps the contents of ppt/theme/theme1.xml (Office Theme) and
# ppt/theme/theme2.xml (Integral) - i.e. the design applied to the deck
# changes from "Integral" colours to the stock "Office Theme" colours.
#
# theme2.xml is the theme actually wired to the slide master (and therefore
# to every slide in the deck), so it is exposed by the PowerPoint object
# model as Presentation.Designs(1).SlideMaster.Theme. We flip its 12 theme
# colours from the Integral palette to the Office Theme palette - this is
# the visible, substantive effect of the commit.

function Get-RGBValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.Designs.Item(1).SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $colorScheme.Item($i + 1).RGB = Get-RGBValue $officeThemeColors[$i]
}
